$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row labels (A2:A10) - replaced with No.1..No.9
$ws.Range("A2").Value = "No.1"
$ws.Range("A3").Value = "No.2"
$ws.Range("A4").Value = "No.3"
$ws.Range("A5").Value = "No.4"
$ws.Range("A6").Value = "No.5"
$ws.Range("A7").Value = "No.6"
$ws.Range("A8").Value = "No.7"
$ws.Range("A9").Value = "No.8"
$ws.Range("A10").Value = "No.9"

# Update header row (B1:F1) - subject names translated to English
$ws.Range("B1").Value = "Japanese"
$ws.Range("C1").Value = "Math"
$ws.Range("D1").Value = "Science"
$ws.Range("E1").Value = "Society"
$ws.Range("F1").Value = "English"

# Update a couple of score values that changed
$ws.Range("F3").Value = 95
$ws.Range("F6").Value = 65

# Update the active selection to match the saved view state
$ws.Range("F4").Select()
